# Atualizado por script em 01-12-2023 20:45
#
# This script:
#  1) Reorders the betting-odds columns (F:V) across several blocks of rows
#     (27-31, 53-60, 83-88, 95-98) to match the freshly re-scraped ordering.
#     Columns A-E (Indice/pais/torneio/temporada/data_partida) stay anchored
#     to their row, only the match/odds columns F:V are permuted.
#  2) Appends two brand-new matches as rows 111 and 112.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Reorder F:V across the affected row blocks.
#    mapping[destinationRow] = sourceRow  (value BEFORE any writes happen)
# ---------------------------------------------------------------------
$mapping = @{
    27 = 29; 28 = 30; 29 = 31; 30 = 28; 31 = 27;
    53 = 54; 54 = 55; 55 = 53;
    57 = 58; 58 = 59; 59 = 60; 60 = 57;
    83 = 86; 85 = 83; 86 = 85; 87 = 88; 88 = 87;
    95 = 97; 96 = 98; 97 = 95; 98 = 96;
}

# Snapshot F:V for every row referenced above (as source or destination)
# BEFORE writing anything, so overlapping swaps don't clobber each other.
$affectedRows = $mapping.Keys + $mapping.Values | Sort-Object -Unique

$snapshot = @{}
foreach ($r in $affectedRows) {
    $snapshot[$r] = $ws.Range("F" + $r + ":V" + $r).Value2
}

foreach ($r in ($mapping.Keys | Sort-Object)) {
    $srcRow = $mapping[$r]
    $ws.Range("F" + $r + ":V" + $r).Value2 = $snapshot[$srcRow]
}

# ---------------------------------------------------------------------
# 2) Append the two newly scraped matches as rows 111 and 112.
#    Copy row 110's formatting down first so the new rows inherit the
#    same cell styles (index/date number formats), then overwrite values.
# ---------------------------------------------------------------------
$ws.Range("A110:V110").Copy($ws.Range("A111:V112"))

$ws.Cells.Item(111, 1).Value = 110
$ws.Cells.Item(111, 2).Value = "portugal"
$ws.Cells.Item(111, 3).Value = "liga-3"
$ws.Cells.Item(111, 4).Value = "2023-2024"
$ws.Cells.Item(111, 5).Value = 45261.66666666666
$ws.Cells.Item(111, 6).Value = "Atletico CP"
$ws.Cells.Item(111, 7).Value = 1
$ws.Cells.Item(111, 8).Value = "Covilha"
$ws.Cells.Item(111, 9).Value = 1
$ws.Cells.Item(111, 10).Value = 2.32
$ws.Cells.Item(111, 11).Value = "24/11/2023 21:42"
$ws.Cells.Item(111, 12).Value = 2.38
$ws.Cells.Item(111, 13).Value = "01/12/2023 15:56"
$ws.Cells.Item(111, 14).Value = 3.21
$ws.Cells.Item(111, 15).Value = "24/11/2023 21:42"
$ws.Cells.Item(111, 16).Value = 3.06
$ws.Cells.Item(111, 17).Value = "01/12/2023 15:56"
$ws.Cells.Item(111, 18).Value = 3.21
$ws.Cells.Item(111, 19).Value = "24/11/2023 21:42"
$ws.Cells.Item(111, 20).Value = 3.36
$ws.Cells.Item(111, 21).Value = "01/12/2023 15:56"
$ws.Cells.Item(111, 22).Value = "https://www.betexplorer.com/football/portugal/liga-3/atletico-cp-covilha/hdftaxea/"

$ws.Cells.Item(112, 1).Value = 111
$ws.Cells.Item(112, 2).Value = "portugal"
$ws.Cells.Item(112, 3).Value = "liga-3"
$ws.Cells.Item(112, 4).Value = "2023-2024"
$ws.Cells.Item(112, 5).Value = 45261.77083333334
$ws.Cells.Item(112, 6).Value = "Alverca"
$ws.Cells.Item(112, 7).Value = 1
$ws.Cells.Item(112, 8).Value = "Caldas"
$ws.Cells.Item(112, 9).Value = 1
$ws.Cells.Item(112, 10).Value = 2.08
$ws.Cells.Item(112, 11).Value = "24/11/2023 21:42"
$ws.Cells.Item(112, 12).Value = 2.26
$ws.Cells.Item(112, 13).Value = "01/12/2023 18:28"
$ws.Cells.Item(112, 14).Value = 3.3
$ws.Cells.Item(112, 15).Value = "24/11/2023 21:42"
$ws.Cells.Item(112, 16).Value = 3.15
$ws.Cells.Item(112, 17).Value = "01/12/2023 18:28"
$ws.Cells.Item(112, 18).Value = 3.7
$ws.Cells.Item(112, 19).Value = "24/11/2023 21:42"
$ws.Cells.Item(112, 20).Value = 3.5
$ws.Cells.Item(112, 21).Value = "01/12/2023 18:25"
$ws.Cells.Item(112, 22).Value = "https://www.betexplorer.com/football/portugal/liga-3/alverca-caldas-sc/OG6T1vBt/"

